# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (col I) and DialogAct
# (col J) values for a set of rows in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I = DAMSLTag, Column J = DialogAct
$updates = @(
    @{ Row = 7;   Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 18;  Tag = "aa"; Label = "Agree/Accept" },
    @{ Row = 19;  Tag = "aa"; Label = "Agree/Accept" },
    @{ Row = 41;  Tag = "ba"; Label = "Appreciation" },
    @{ Row = 42;  Tag = "b";  Label = "Acknowledge (Backchannel)" },
    @{ Row = 43;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 46;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 54;  Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 62;  Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 71;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 75;  Tag = "ba"; Label = "Appreciation" },
    @{ Row = 79;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 104; Tag = "sv"; Label = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Label
}
